$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K24").Value = -3.87
$ws.Range("J25").Value = -3.86
$ws.Range("K25").Value = -0.5599999999999999
$ws.Range("I26").Value = -3.94
$ws.Range("J26").Value = -0.64
$ws.Range("K26").Value = 0.5599999999999999
$ws.Range("H27").Value = -3.97
$ws.Range("I27").Value = -0.67
$ws.Range("J27").Value = 0.5299999999999999
$ws.Range("K27").Value = 0.5299999999999999
$ws.Range("G28").Value = -3.91
$ws.Range("H28").Value = -0.61
$ws.Range("I28").Value = 0.59
$ws.Range("J28").Value = 0.59
$ws.Range("K28").Value = 0.19
$ws.Range("F29").Value = -3.93
$ws.Range("G29").Value = -0.63
$ws.Range("H29").Value = 0.57
$ws.Range("I29").Value = 0.57
$ws.Range("J29").Value = 0.17
$ws.Range("K29").Value = -0.23
$ws.Range("E30").Value = -3.95
$ws.Range("F30").Value = -0.65
$ws.Range("G30").Value = 0.5499999999999999
$ws.Range("H30").Value = 0.5499999999999999
$ws.Range("I30").Value = 0.15
$ws.Range("J30").Value = -0.25
$ws.Range("K30").Value = 1.35
$ws.Range("D31").Value = -3.94
$ws.Range("E31").Value = -0.64
$ws.Range("F31").Value = 0.5599999999999999
$ws.Range("G31").Value = 0.5599999999999999
$ws.Range("H31").Value = 0.16
$ws.Range("I31").Value = -0.24
$ws.Range("J31").Value = 1.36
$ws.Range("K31").Value = 0.5599999999999999
$ws.Range("C32").Value = -4.05
$ws.Range("D32").Value = -0.75
$ws.Range("E32").Value = 0.45
$ws.Range("F32").Value = 0.45
$ws.Range("G32").Value = 0.05000000000000002
$ws.Range("H32").Value = -0.35
$ws.Range("I32").Value = 1.25
$ws.Range("J32").Value = 0.45
$ws.Range("K32").Value = -0.04999999999999999
$ws.Range("B33").Value = -4.38
$ws.Range("C33").Value = -1.08
$ws.Range("D33").Value = 0.12
$ws.Range("E33").Value = 0.12
$ws.Range("F33").Value = -0.28
$ws.Range("G33").Value = -0.6799999999999999
$ws.Range("H33").Value = 0.9199999999999999
$ws.Range("I33").Value = 0.12
$ws.Range("J33").Value = -0.38
$ws.Range("K33").Value = -0.58
$ws.Range("B34").Value = -0.7
$ws.Range("C34").Value = 0.5
$ws.Range("D34").Value = 0.5
$ws.Range("E34").Value = 0.1
$ws.Range("F34").Value = -0.3
$ws.Range("G34").Value = 1.3
$ws.Range("H34").Value = 0.5
$ws.Range("I34").Value = 0.0
$ws.Range("J34").Value = -0.2
$ws.Range("K34").Value = -0.7999999999999999
$ws.Range("B35").Value = 0.61
$ws.Range("C35").Value = 0.61
$ws.Range("D35").Value = 0.21
$ws.Range("E35").Value = -0.19
$ws.Range("F35").Value = 1.41
$ws.Range("G35").Value = 0.61
$ws.Range("H35").Value = 0.11
$ws.Range("I35").Value = -0.08999999999999997
$ws.Range("J35").Value = -0.69
$ws.Range("K35").Value = 0.61
$ws.Range("B36").Value = 0.63
$ws.Range("C36").Value = 0.23
$ws.Range("D36").Value = -0.17
$ws.Range("E36").Value = 1.43
$ws.Range("F36").Value = 0.63
$ws.Range("G36").Value = 0.13
$ws.Range("H36").Value = -0.07000000000000002
$ws.Range("I36").Value = -0.6699999999999999
$ws.Range("J36").Value = 0.63
$ws.Range("K36").Value = -0.17
$ws.Range("B37").Value = 0.15
$ws.Range("C37").Value = -0.25
$ws.Range("D37").Value = 1.35
$ws.Range("E37").Value = 0.55
$ws.Range("F37").Value = 0.05000000000000002
$ws.Range("G37").Value = -0.15
$ws.Range("H37").Value = -0.75
$ws.Range("I37").Value = 0.55
$ws.Range("J37").Value = -0.25
$ws.Range("K37").Value = -0.04999999999999999
$ws.Range("B38").Value = -0.28
$ws.Range("C38").Value = 1.32
$ws.Range("D38").Value = 0.5199999999999999
$ws.Range("E38").Value = 0.01999999999999998
$ws.Range("F38").Value = -0.18
$ws.Range("G38").Value = -0.78
$ws.Range("H38").Value = 0.5199999999999999
$ws.Range("I38").Value = -0.28
$ws.Range("J38").Value = -0.08000000000000003
$ws.Range("K38").Value = -0.38
$ws.Range("B39").Value = 1.45
$ws.Range("C39").Value = 0.6499999999999999
$ws.Range("D39").Value = 0.15
$ws.Range("E39").Value = -0.05000000000000002
$ws.Range("F39").Value = -0.65
$ws.Range("G39").Value = 0.6499999999999999
$ws.Range("H39").Value = -0.15
$ws.Range("I39").Value = 0.04999999999999998
$ws.Range("J39").Value = -0.25
$ws.Range("K39").Value = 0.15
$ws.Range("B40").Value = 0.41
$ws.Range("C40").Value = -0.09
$ws.Range("D40").Value = -0.29
$ws.Range("E40").Value = -0.8899999999999999
$ws.Range("F40").Value = 0.41
$ws.Range("G40").Value = -0.39
$ws.Range("H40").Value = -0.19
$ws.Range("I40").Value = -0.49
$ws.Range("J40").Value = -0.09
$ws.Range("K40").Value = 0.71
$ws.Range("B41").Value = -0.01999999999999999
$ws.Range("C41").Value = -0.22
$ws.Range("D41").Value = -0.82
$ws.Range("E41").Value = 0.48
$ws.Range("F41").Value = -0.32
$ws.Range("G41").Value = -0.12
$ws.Range("H41").Value = -0.42
$ws.Range("I41").Value = -0.01999999999999999
$ws.Range("J41").Value = 0.78
$ws.Range("K41").Value = 1.08
$ws.Range("B42").Value = -0.09000000000000004
$ws.Range("C42").Value = -0.69
$ws.Range("D42").Value = 0.61
$ws.Range("E42").Value = -0.1900000000000001
$ws.Range("F42").Value = 0.009999999999999962
$ws.Range("G42").Value = -0.29
$ws.Range("H42").Value = 0.11
$ws.Range("I42").Value = 0.91
$ws.Range("J42").Value = 1.21
$ws.Range("K42").Value = -0.89
$ws.Range("B43").Value = -0.41
$ws.Range("C43").Value = 0.8899999999999999
$ws.Range("D43").Value = 0.08999999999999997
$ws.Range("E43").Value = 0.29
$ws.Range("F43").Value = -0.01000000000000001
$ws.Range("G43").Value = 0.39
$ws.Range("H43").Value = 1.19
$ws.Range("I43").Value = 1.49
$ws.Range("J43").Value = -0.6100000000000001
$ws.Range("K43").Value = 0.99
$ws.Range("B44").Value = 0.57
$ws.Range("C44").Value = -0.23
$ws.Range("D44").Value = -0.03000000000000001
$ws.Range("E44").Value = -0.33
$ws.Range("F44").Value = 0.06999999999999999
$ws.Range("G44").Value = 0.87
$ws.Range("H44").Value = 1.17
$ws.Range("I44").Value = -0.93
$ws.Range("J44").Value = 0.6699999999999999
$ws.Range("B45").Value = -0.34
$ws.Range("C45").Value = -0.14
$ws.Range("D45").Value = -0.44
$ws.Range("E45").Value = -0.04000000000000001
$ws.Range("F45").Value = 0.76
$ws.Range("G45").Value = 1.06
$ws.Range("H45").Value = -1.04
$ws.Range("I45").Value = 0.5599999999999999
$ws.Range("B46").Value = -0.01000000000000001
$ws.Range("C46").Value = -0.31
$ws.Range("D46").Value = 0.09
$ws.Range("E46").Value = 0.89
$ws.Range("F46").Value = 1.19
$ws.Range("G46").Value = -0.91
$ws.Range("H46").Value = 0.69
$ws.Range("B47").Value = -0.27
$ws.Range("C47").Value = 0.13
$ws.Range("D47").Value = 0.93
$ws.Range("E47").Value = 1.23
$ws.Range("F47").Value = -0.87
$ws.Range("G47").Value = 0.73
$ws.Range("B48").Value = 0.19
$ws.Range("C48").Value = 0.99
$ws.Range("D48").Value = 1.29
$ws.Range("E48").Value = -0.81
$ws.Range("F48").Value = 0.7899999999999999
$ws.Range("B49").Value = 0.9600000000000001
$ws.Range("C49").Value = 1.26
$ws.Range("D49").Value = -0.84
$ws.Range("E49").Value = 0.76
$ws.Range("B50").Value = 1.14
$ws.Range("C50").Value = -0.9600000000000001
$ws.Range("D50").Value = 0.6399999999999999
$ws.Range("B51").Value = -1.03
$ws.Range("C51").Value = 0.57
$ws.Range("B52").Value = 0.71
